$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "test"
$ws.Range("B4").Value = "1pmk193ci/a"
$ws.Range("C4").Value = "Rack A"
$ws.Range("D4").Value = 10
